$wb = $excel.ActiveWorkbook

# --- ALERTS: rows 4-5 ---
$ws = $wb.Worksheets.Item('ALERTS')
$arr = New-Object 'object[,]' 2,6
$arr[0,0] = "'2026-01-30"
$arr[0,1] = '13:06:14'
$arr[0,2] = '13:00'
$arr[0,3] = 'Bathroom'
$arr[0,4] = 'MINIMAL'
$arr[0,5] = 'MINIMAL ALERT: Bathroom occupied, no motion > 20s.'
$arr[1,0] = "'2026-01-30"
$arr[1,1] = '13:06:45'
$arr[1,2] = '13:00'
$arr[1,3] = 'Bathroom'
$arr[1,4] = 'MINIMAL'
$arr[1,5] = 'MINIMAL ALERT: Bathroom occupied, no motion > 20s.'
$ws.Range('A4:F5').Value = $arr

# --- PIR: rows 75-96 ---
$ws = $wb.Worksheets.Item('PIR')
$arr = New-Object 'object[,]' 22,6
$arr[0,0] = "'2026-01-30"
$arr[0,1] = '13:03:08'
$arr[0,2] = '13:00'
$arr[0,3] = 'Bathroom'
$arr[0,4] = 'No Motion'
$arr[0,5] = 'Inactive'
$arr[1,0] = "'2026-01-30"
$arr[1,1] = '13:03:08'
$arr[1,2] = '13:00'
$arr[1,3] = 'Bathroom'
$arr[1,4] = 'No Motion'
$arr[1,5] = 'Inactive'
$arr[2,0] = "'2026-01-30"
$arr[2,1] = '13:03:13'
$arr[2,2] = '13:00'
$arr[2,3] = 'Bathroom'
$arr[2,4] = 'No Motion'
$arr[2,5] = 'Inactive'
$arr[3,0] = "'2026-01-30"
$arr[3,1] = '13:03:18'
$arr[3,2] = '13:00'
$arr[3,3] = 'Bathroom'
$arr[3,4] = 'No Motion'
$arr[3,5] = 'Inactive'
$arr[4,0] = "'2026-01-30"
$arr[4,1] = '13:03:23'
$arr[4,2] = '13:00'
$arr[4,3] = 'Bathroom'
$arr[4,4] = 'No Motion'
$arr[4,5] = 'Inactive'
$arr[5,0] = "'2026-01-30"
$arr[5,1] = '13:03:28'
$arr[5,2] = '13:00'
$arr[5,3] = 'Bathroom'
$arr[5,4] = 'No Motion'
$arr[5,5] = 'Inactive'
$arr[6,0] = "'2026-01-30"
$arr[6,1] = '13:03:28'
$arr[6,2] = '13:00'
$arr[6,3] = 'Living Room'
$arr[6,4] = 'RECOVERY_DETECTION'
$arr[6,5] = 'Inactive'
$arr[7,0] = "'2026-01-30"
$arr[7,1] = '13:03:33'
$arr[7,2] = '13:00'
$arr[7,3] = 'Bathroom'
$arr[7,4] = 'No Motion'
$arr[7,5] = 'Inactive'
$arr[8,0] = "'2026-01-30"
$arr[8,1] = '13:03:38'
$arr[8,2] = '13:00'
$arr[8,3] = 'Bathroom'
$arr[8,4] = 'No Motion'
$arr[8,5] = 'Inactive'
$arr[9,0] = "'2026-01-30"
$arr[9,1] = '13:03:43'
$arr[9,2] = '13:00'
$arr[9,3] = 'Bathroom'
$arr[9,4] = 'No Motion'
$arr[9,5] = 'Inactive'
$arr[10,0] = "'2026-01-30"
$arr[10,1] = '13:05:49'
$arr[10,2] = '13:00'
$arr[10,3] = 'Bathroom'
$arr[10,4] = 'No Motion'
$arr[10,5] = 'Inactive'
$arr[11,0] = "'2026-01-30"
$arr[11,1] = '13:05:54'
$arr[11,2] = '13:00'
$arr[11,3] = 'Bathroom'
$arr[11,4] = 'No Motion'
$arr[11,5] = 'Inactive'
$arr[12,0] = "'2026-01-30"
$arr[12,1] = '13:05:59'
$arr[12,2] = '13:00'
$arr[12,3] = 'Bathroom'
$arr[12,4] = 'No Motion'
$arr[12,5] = 'Inactive'
$arr[13,0] = "'2026-01-30"
$arr[13,1] = '13:06:04'
$arr[13,2] = '13:00'
$arr[13,3] = 'Bathroom'
$arr[13,4] = 'No Motion'
$arr[13,5] = 'Inactive'
$arr[14,0] = "'2026-01-30"
$arr[14,1] = '13:06:09'
$arr[14,2] = '13:00'
$arr[14,3] = 'Bathroom'
$arr[14,4] = 'No Motion'
$arr[14,5] = 'Inactive'
$arr[15,0] = "'2026-01-30"
$arr[15,1] = '13:06:14'
$arr[15,2] = '13:00'
$arr[15,3] = 'Bathroom'
$arr[15,4] = 'No Motion'
$arr[15,5] = 'Inactive'
$arr[16,0] = "'2026-01-30"
$arr[16,1] = '13:06:19'
$arr[16,2] = '13:00'
$arr[16,3] = 'Bathroom'
$arr[16,4] = 'No Motion'
$arr[16,5] = 'Inactive'
$arr[17,0] = "'2026-01-30"
$arr[17,1] = '13:06:24'
$arr[17,2] = '13:00'
$arr[17,3] = 'Bathroom'
$arr[17,4] = 'No Motion'
$arr[17,5] = 'Inactive'
$arr[18,0] = "'2026-01-30"
$arr[18,1] = '13:06:29'
$arr[18,2] = '13:00'
$arr[18,3] = 'Bathroom'
$arr[18,4] = 'No Motion'
$arr[18,5] = 'Inactive'
$arr[19,0] = "'2026-01-30"
$arr[19,1] = '13:06:34'
$arr[19,2] = '13:00'
$arr[19,3] = 'Bathroom'
$arr[19,4] = 'No Motion'
$arr[19,5] = 'Inactive'
$arr[20,0] = "'2026-01-30"
$arr[20,1] = '13:06:39'
$arr[20,2] = '13:00'
$arr[20,3] = 'Bathroom'
$arr[20,4] = 'No Motion'
$arr[20,5] = 'Inactive'
$arr[21,0] = "'2026-01-30"
$arr[21,1] = '13:06:44'
$arr[21,2] = '13:00'
$arr[21,3] = 'Bathroom'
$arr[21,4] = 'No Motion'
$arr[21,5] = 'Inactive'
$ws.Range('A75:F96').Value = $arr

# --- Humidity: rows 56-56 ---
$ws = $wb.Worksheets.Item('Humidity')
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = "'2026-01-30"
$arr[0,1] = '13:03:07'
$arr[0,2] = '13:00'
$arr[0,3] = 'Bathroom'
$arr[0,4] = "'87.3%"
$arr[0,5] = 'Active'
$ws.Range('A56:F56').Value = $arr

# --- Temperature: rows 56-56 ---
$ws = $wb.Worksheets.Item('Temperature')
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = "'2026-01-30"
$arr[0,1] = '13:03:07'
$arr[0,2] = '13:00'
$arr[0,3] = 'Bathroom'
$arr[0,4] = '22.6C'
$arr[0,5] = 'Active'
$ws.Range('A56:F56').Value = $arr

# --- Proximity: rows 33-42 ---
$ws = $wb.Worksheets.Item('Proximity')
$arr = New-Object 'object[,]' 10,6
$arr[0,0] = "'2026-01-30"
$arr[0,1] = '13:03:07'
$arr[0,2] = '13:00'
$arr[0,3] = 'Bathroom Door'
$arr[0,4] = 'EXIT'
$arr[0,5] = 'User EXITED Bathroom'
$arr[1,0] = "'2026-01-30"
$arr[1,1] = '13:03:09'
$arr[1,2] = '13:00'
$arr[1,3] = 'Bathroom Door'
$arr[1,4] = 'ENTER'
$arr[1,5] = 'User ENTERED Bathroom'
$arr[2,0] = "'2026-01-30"
$arr[2,1] = '13:03:11'
$arr[2,2] = '13:00'
$arr[2,3] = 'Bathroom Door'
$arr[2,4] = 'EXIT'
$arr[2,5] = 'User EXITED Bathroom'
$arr[3,0] = "'2026-01-30"
$arr[3,1] = '13:03:15'
$arr[3,2] = '13:00'
$arr[3,3] = 'Bathroom Door'
$arr[3,4] = 'ENTER'
$arr[3,5] = 'User ENTERED Bathroom'
$arr[4,0] = "'2026-01-30"
$arr[4,1] = '13:03:22'
$arr[4,2] = '13:00'
$arr[4,3] = 'Bathroom Door'
$arr[4,4] = 'EXIT'
$arr[4,5] = 'User EXITED Bathroom'
$arr[5,0] = "'2026-01-30"
$arr[5,1] = '13:03:30'
$arr[5,2] = '13:00'
$arr[5,3] = 'Bathroom Door'
$arr[5,4] = 'ENTER'
$arr[5,5] = 'User ENTERED Bathroom'
$arr[6,0] = "'2026-01-30"
$arr[6,1] = '13:03:37'
$arr[6,2] = '13:00'
$arr[6,3] = 'Bathroom Door'
$arr[6,4] = 'EXIT'
$arr[6,5] = 'User EXITED Bathroom'
$arr[7,0] = "'2026-01-30"
$arr[7,1] = '13:05:50'
$arr[7,2] = '13:00'
$arr[7,3] = 'Bathroom Door'
$arr[7,4] = 'ENTER'
$arr[7,5] = 'User ENTERED Bathroom'
$arr[8,0] = "'2026-01-30"
$arr[8,1] = '13:06:11'
$arr[8,2] = '13:00'
$arr[8,3] = 'Bathroom Door'
$arr[8,4] = 'EXIT'
$arr[8,5] = 'User EXITED Bathroom'
$arr[9,0] = "'2026-01-30"
$arr[9,1] = '13:06:21'
$arr[9,2] = '13:00'
$arr[9,3] = 'Bathroom Door'
$arr[9,4] = 'ENTER'
$arr[9,5] = 'User ENTERED Bathroom'
$ws.Range('A33:F42').Value = $arr

# --- mmWave: rows 31-38 ---
$ws = $wb.Worksheets.Item('mmWave')
$arr = New-Object 'object[,]' 8,6
$arr[0,0] = "'2026-01-30"
$arr[0,1] = '13:03:08'
$arr[0,2] = '13:00'
$arr[0,3] = 'Living Room'
$arr[0,4] = 'PRESENCE_DETECTED'
$arr[0,5] = 'Active'
$arr[1,0] = "'2026-01-30"
$arr[1,1] = '13:03:28'
$arr[1,2] = '13:00'
$arr[1,3] = 'Living Room'
$arr[1,4] = 'FALL_DETECTED'
$arr[1,5] = 'EMERGENCY'
$arr[2,0] = "'2026-01-30"
$arr[2,1] = '13:03:29'
$arr[2,2] = '13:00'
$arr[2,3] = 'Living Room'
$arr[2,4] = 'PRESENCE_DETECTED'
$arr[2,5] = 'Active'
$arr[3,0] = "'2026-01-30"
$arr[3,1] = '13:03:39'
$arr[3,2] = '13:00'
$arr[3,3] = 'Living Room'
$arr[3,4] = 'PRESENCE_DETECTED'
$arr[3,5] = 'Active'
$arr[4,0] = "'2026-01-30"
$arr[4,1] = '13:05:48'
$arr[4,2] = '13:00'
$arr[4,3] = 'Living Room'
$arr[4,4] = 'FALL_DETECTED'
$arr[4,5] = 'EMERGENCY'
$arr[5,0] = "'2026-01-30"
$arr[5,1] = '13:05:48'
$arr[5,2] = '13:00'
$arr[5,3] = 'Living Room'
$arr[5,4] = 'FALL_DETECTED'
$arr[5,5] = 'EMERGENCY'
$arr[6,0] = "'2026-01-30"
$arr[6,1] = '13:05:57'
$arr[6,2] = '13:00'
$arr[6,3] = 'Living Room'
$arr[6,4] = 'PRESENCE_DETECTED'
$arr[6,5] = 'Active'
$arr[7,0] = "'2026-01-30"
$arr[7,1] = '13:06:19'
$arr[7,2] = '13:00'
$arr[7,3] = 'Living Room'
$arr[7,4] = 'FALL_DETECTED'
$arr[7,5] = 'EMERGENCY'
$ws.Range('A31:F38').Value = $arr
